# Update Leve profit-tracking figures across sheets (scheduled data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 7124.875
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 7799.8
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 7799.8
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -8295.799999999999
# Row 67
$ws.Range("H67").Value = 7124.875
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 7799.8
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 7799.8
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -9515.799999999999
# Row 92
$ws.Range("H92").Value = 2815.75
$ws.Range("I92").Value = 653.4
$ws.Range("J92").Value = 6419.6665
$ws.Range("K92").Value = 653.4
$ws.Range("L92").Value = 6419.6665
$ws.Range("M92").Value = 594.6
$ws.Range("N92").Value = -8915.666499999999
# Row 111
$ws.Range("H111").Value = 778.4
$ws.Range("I111").Value = 573
$ws.Range("J111").Value = 1600
$ws.Range("K111").Value = 1719
$ws.Range("L111").Value = 4800
$ws.Range("M111").Value = 1348
$ws.Range("N111").Value = -10934
# Row 116
$ws.Range("H116").Value = 10126.25
$ws.Range("I116").Value = 9501
$ws.Range("J116").Value = 10751.5
$ws.Range("K116").Value = 9501
$ws.Range("L116").Value = 10751.5
$ws.Range("M116").Value = -6059
$ws.Range("N116").Value = -17635.5

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 224.41667
$ws.Range("I97").Value = 228.1
$ws.Range("J97").Value = 206
$ws.Range("K97").Value = 228.1
$ws.Range("L97").Value = 206
$ws.Range("M97").Value = 267.9
$ws.Range("N97").Value = -1198

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1849.7
$ws.Range("I5").Value = 183
$ws.Range("J5").Value = 2564
$ws.Range("K5").Value = 183
$ws.Range("L5").Value = 2564
$ws.Range("M5").Value = -70
$ws.Range("N5").Value = -2790

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3822.5
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 5645
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 5645
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -6219
# Row 58
$ws.Range("H58").Value = 3917.1
$ws.Range("I58").Value = 1163
$ws.Range("J58").Value = 9031.857
$ws.Range("K58").Value = 1163
$ws.Range("L58").Value = 9031.857
$ws.Range("M58").Value = -960
$ws.Range("N58").Value = -9437.857
# Row 113
$ws.Range("H113").Value = 3822.5
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 5645
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 5645
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -9985
# Row 132
$ws.Range("H132").Value = 4308.3
$ws.Range("I132").Value = 3930.875
$ws.Range("J132").Value = 5818
$ws.Range("K132").Value = 11792.625
$ws.Range("L132").Value = 17454
$ws.Range("M132").Value = -9262.625
$ws.Range("N132").Value = -22514
# Row 136
$ws.Range("H136").Value = 3917.1
$ws.Range("I136").Value = 1163
$ws.Range("J136").Value = 9031.857
$ws.Range("K136").Value = 3489
$ws.Range("L136").Value = 27095.571
$ws.Range("M136").Value = -939
$ws.Range("N136").Value = -32195.571

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 125462.125
$ws.Range("I11").Value = 200731.4
$ws.Range("J11").Value = 13.333333
$ws.Range("K11").Value = 602194.2
$ws.Range("L11").Value = 39.999999
$ws.Range("M11").Value = -602054.2
$ws.Range("N11").Value = -319.999999
# Row 114
$ws.Range("H114").Value = 70.42856999999999
$ws.Range("I114").Value = 250
$ws.Range("J114").Value = 40.5
$ws.Range("K114").Value = 750
$ws.Range("L114").Value = 121.5
$ws.Range("M114").Value = 2504
$ws.Range("N114").Value = -6629.5

$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 34999
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 34999
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 34999
$ws.Range("N26").Value = -35559
# Row 50
$ws.Range("H50").Value = 34999
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 34999
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 34999
$ws.Range("N50").Value = -35995
# Row 80
$ws.Range("H80").Value = 6334.9287
$ws.Range("I80").Value = 6586.625
$ws.Range("J80").Value = 5999.3335
$ws.Range("K80").Value = 6586.625
$ws.Range("L80").Value = 5999.3335
$ws.Range("M80").Value = -5588.625
$ws.Range("N80").Value = -7995.3335
# Row 83
$ws.Range("H83").Value = 6334.9287
$ws.Range("I83").Value = 6586.625
$ws.Range("J83").Value = 5999.3335
$ws.Range("K83").Value = 32933.125
$ws.Range("L83").Value = 29996.6675
$ws.Range("M83").Value = -27941.125
$ws.Range("N83").Value = -39980.6675
# Row 107
$ws.Range("H107").Value = 1054.0869
$ws.Range("I107").Value = 462.72726
$ws.Range("J107").Value = 1596.1666
$ws.Range("K107").Value = 462.72726
$ws.Range("L107").Value = 1596.1666
$ws.Range("M107").Value = 1457.27274
$ws.Range("N107").Value = -5436.1666
# Row 126
$ws.Range("H126").Value = 7194.6
$ws.Range("I126").Value = 4319.6665
$ws.Range("J126").Value = 11507
$ws.Range("K126").Value = 12958.9995
$ws.Range("L126").Value = 34521
$ws.Range("M126").Value = -10488.9995
$ws.Range("N126").Value = -39461

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 60994.75
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 60994.75
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 60994.75
$ws.Range("N6").Value = -61218.75
# Row 40
$ws.Range("H40").Value = 7982.4165
$ws.Range("I40").Value = 6764.8887
$ws.Range("J40").Value = 11635
$ws.Range("K40").Value = 6764.8887
$ws.Range("L40").Value = 11635
$ws.Range("M40").Value = -6628.8887
$ws.Range("N40").Value = -11907
# Row 116
$ws.Range("H116").Value = 350000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 350000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 350000
$ws.Range("N116").Value = -359178
# Row 117
$ws.Range("H117").Value = 71842
$ws.Range("I117").Value = 70000
$ws.Range("J117").Value = 73684
$ws.Range("K117").Value = 70000
$ws.Range("L117").Value = 73684
$ws.Range("M117").Value = -65411
$ws.Range("N117").Value = -82862
# Row 118
$ws.Range("H118").Value = 100000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 100000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314
# Row 132
$ws.Range("H132").Value = 3258.4856
$ws.Range("I132").Value = 1734.5667
$ws.Range("J132").Value = 12402
$ws.Range("K132").Value = 5203.7001
$ws.Range("L132").Value = 37206
$ws.Range("M132").Value = -2673.7001
$ws.Range("N132").Value = -42266
# Row 136
$ws.Range("H136").Value = 8700.416999999999
$ws.Range("I136").Value = 2065.8333
$ws.Range("J136").Value = 15335
$ws.Range("K136").Value = 6197.499899999999
$ws.Range("L136").Value = 46005
$ws.Range("M136").Value = -3647.499899999999
$ws.Range("N136").Value = -51105

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5138.143
$ws.Range("I81").Value = 2993
$ws.Range("J81").Value = 10501
$ws.Range("K81").Value = 5986
$ws.Range("L81").Value = 21002
$ws.Range("M81").Value = -4925
$ws.Range("N81").Value = -23124
# Row 84
$ws.Range("H84").Value = 5138.143
$ws.Range("I84").Value = 2993
$ws.Range("J84").Value = 10501
$ws.Range("K84").Value = 29930
$ws.Range("L84").Value = 105010
$ws.Range("M84").Value = -24626
$ws.Range("N84").Value = -115618
# Row 107
$ws.Range("H107").Value = 451.05
$ws.Range("I107").Value = 340.1111
$ws.Range("J107").Value = 1449.5
$ws.Range("K107").Value = 1020.3333
$ws.Range("L107").Value = 4348.5
$ws.Range("M107").Value = 899.6667
$ws.Range("N107").Value = -8188.5
# Row 116
$ws.Range("H116").Value = 151948.4
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 151948.4
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 151948.4
$ws.Range("N116").Value = -161126.4
# Row 118
$ws.Range("H118").Value = 65977
$ws.Range("I118").Value = 78946
$ws.Range("J118").Value = 62734.75
$ws.Range("K118").Value = 78946
$ws.Range("L118").Value = 62734.75
$ws.Range("M118").Value = -77289
$ws.Range("N118").Value = -66048.75
# Row 132
$ws.Range("H132").Value = 8302.5625
$ws.Range("I132").Value = 7589.067
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 22767.201
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -20237.201
$ws.Range("N132").Value = -62075
